$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 48.484375
$ws.Range("H2").Value = 145.453125
$ws.Range("I2").Value = 0.7776469276297807
$ws.Range("J2").Value = 0.7776469276297806
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 7443.758449770832
$ws.Range("R2").Value = 66993.8260479375
$ws.Range("S2").Value = 0.2466857003107477
$ws.Range("T2").Value = 0.2466857003107478

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 48.484375
$ws.Range("H3").Value = 145.453125
$ws.Range("I3").Value = 0.7776469276297807
$ws.Range("J3").Value = 0.7776469276297806
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 8184.148083979167
$ws.Range("R3").Value = 73657.33275581252
$ws.Range("S3").Value = 0.2712221675604507
$ws.Range("T3").Value = 0.2712221675604506

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 48.484375
$ws.Range("H4").Value = 145.453125
$ws.Range("I4").Value = 0.7776469276297807
$ws.Range("J4").Value = 0.7776469276297806
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 3301.316770364583
$ws.Range("R4").Value = 29711.85093328125
$ws.Range("S4").Value = 0.109405436103327
$ws.Range("T4").Value = 0.109405436103327

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 48.484375
$ws.Range("H5").Value = 145.453125
$ws.Range("I5").Value = 0.7776469276297807
$ws.Range("J5").Value = 0.7776469276297806
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 4536.327723734375
$ws.Range("R5").Value = 40826.94951360938
$ws.Range("S5").Value = 0.1503336236552553
$ws.Range("T5").Value = 0.1503336236552552

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.47148733333333
$ws.Range("H6").Value = 34.414462
$ws.Range("I6").Value = 0.1839926137051496
$ws.Range("J6").Value = 0.1839926137051496
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 1761.206177638447
$ws.Range("R6").Value = 15850.85559874602
$ws.Range("S6").Value = 0.05836626514066038
$ws.Range("T6").Value = 0.05836626514066038

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.47148733333333
$ws.Range("H7").Value = 34.414462
$ws.Range("I7").Value = 0.1839926137051496
$ws.Range("J7").Value = 0.1839926137051496
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 1936.3836510111
$ws.Range("R7").Value = 17427.4528590999
$ws.Range("S7").Value = 0.0641716359072159
$ws.Range("T7").Value = 0.06417163590721589

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.47148733333333
$ws.Range("H8").Value = 34.414462
$ws.Range("I8").Value = 0.1839926137051496
$ws.Range("J8").Value = 0.1839926137051496
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 781.0972816409045
$ws.Range("R8").Value = 7029.87553476814
$ws.Range("S8").Value = 0.02588551619892233
$ws.Range("T8").Value = 0.02588551619892233

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.47148733333333
$ws.Range("H9").Value = 34.414462
$ws.Range("I9").Value = 0.1839926137051496
$ws.Range("J9").Value = 0.1839926137051496
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 1073.303018192309
$ws.Range("R9").Value = 9659.727163730779
$ws.Range("S9").Value = 0.03556919645835099
$ws.Range("T9").Value = 0.03556919645835098

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.356432
$ws.Range("H10").Value = 1.069296
$ws.Range("I10").Value = 0.005716857228930723
$ws.Range("J10").Value = 0.005716857228930722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 54.72265470615466
$ws.Range("R10").Value = 492.503892355392
$ws.Range("S10").Value = 0.00181350543413544
$ws.Range("T10").Value = 0.00181350543413544

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.356432
$ws.Range("H11").Value = 1.069296
$ws.Range("I11").Value = 0.005716857228930723
$ws.Range("J11").Value = 0.005716857228930722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 60.16561562088535
$ws.Range("R11").Value = 541.4905405879681
$ws.Range("S11").Value = 0.00199388482635708
$ws.Range("T11").Value = 0.00199388482635708

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.356432
$ws.Range("H12").Value = 1.069296
$ws.Range("I12").Value = 0.005716857228930723
$ws.Range("J12").Value = 0.005716857228930722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 24.26957012634667
$ws.Range("R12").Value = 218.42613113712
$ws.Range("S12").Value = 0.0008042920714391193
$ws.Range("T12").Value = 0.0008042920714391193

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.356432
$ws.Range("H13").Value = 1.069296
$ws.Range("I13").Value = 0.005716857228930723
$ws.Range("J13").Value = 0.005716857228930722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 33.348730662736
$ws.Range("R13").Value = 300.138575964624
$ws.Range("S13").Value = 0.001105174896999084
$ws.Range("T13").Value = 0.001105174896999084

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.035248333333333
$ws.Range("H14").Value = 6.105745
$ws.Range("I14").Value = 0.03264360143613892
$ws.Range("J14").Value = 0.03264360143613892
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 312.469676645971
$ws.Range("R14").Value = 2812.22708981374
$ws.Range("S14").Value = 0.01035522599630531
$ws.Range("T14").Value = 0.01035522599630532

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.035248333333333
$ws.Range("H15").Value = 6.105745
$ws.Range("I15").Value = 0.03264360143613892
$ws.Range("J15").Value = 0.03264360143613892
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 343.5493135194956
$ws.Range("R15").Value = 3091.94382167546
$ws.Range("S15").Value = 0.01138520326374139
$ws.Range("T15").Value = 0.01138520326374139

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.035248333333333
$ws.Range("H16").Value = 6.105745
$ws.Range("I16").Value = 0.03264360143613892
$ws.Range("J16").Value = 0.03264360143613892
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 138.5807170802944
$ws.Range("R16").Value = 1247.22645372265
$ws.Range("S16").Value = 0.004592556498601926
$ws.Range("T16").Value = 0.004592556498601926

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.035248333333333
$ws.Range("H17").Value = 6.105745
$ws.Range("I17").Value = 0.03264360143613892
$ws.Range("J17").Value = 0.03264360143613892
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 190.4232742854616
$ws.Range("R17").Value = 1713.809468569155
$ws.Range("S17").Value = 0.0063106156774903
$ws.Range("T17").Value = 0.0063106156774903
